$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update title to new rich-text string (keeps existing s="2" style on B2)
$ws.Cells.Item(2,2).Value = "Project Bloom\n<size=30>A Quest for Home</size>"

# Rows 11-73: reorder Key/Value pairs (existing per-row styles on column B are untouched)
$ws.Cells.Item(11,1).Value = "continue"
$ws.Cells.Item(11,2).Value = "CONTINUE"
$ws.Cells.Item(12,1).Value = "newGame"
$ws.Cells.Item(12,2).Value = "NEW GAME"
$ws.Cells.Item(13,1).Value = "criteria"
$ws.Cells.Item(13,2).Value = "CRITERIA"
$ws.Cells.Item(14,1).Value = "hotspotAnalyze_title"
$ws.Cells.Item(14,2).Value = "Atmospheric Reading"
$ws.Cells.Item(15,1).Value = "analyzing"
$ws.Cells.Item(15,2).Value = "ANALYZING"
$ws.Cells.Item(16,1).Value = "incompatible"
$ws.Cells.Item(16,2).Value = "INCOMPATIBLE!"
$ws.Cells.Item(17,1).Value = "investigate"
$ws.Cells.Item(17,2).Value = "INVESTIGATE"
$ws.Cells.Item(18,1).Value = "launch"
$ws.Cells.Item(18,2).Value = "LAUNCH"
$ws.Cells.Item(19,1).Value = "back"
$ws.Cells.Item(19,2).Value = "BACK"
$ws.Cells.Item(20,1).Value = "cancel"
$ws.Cells.Item(20,2).Value = "CANCEL"
$ws.Cells.Item(21,1).Value = "new_house"
$ws.Cells.Item(21,2).Value = "A new house is available! Deploy it to progress."
$ws.Cells.Item(22,1).Value = "victory"
$ws.Cells.Item(22,2).Value = "VICTORY"
$ws.Cells.Item(23,1).Value = "population"
$ws.Cells.Item(23,2).Value = "Population"
$ws.Cells.Item(24,1).Value = "houses_deployed"
$ws.Cells.Item(24,2).Value = "Homes Deployed"
$ws.Cells.Item(25,1).Value = "season_title"
$ws.Cells.Item(25,2).Value = "Season"
$ws.Cells.Item(26,1).Value = "season_winter"
$ws.Cells.Item(26,2).Value = "Winter"
$ws.Cells.Item(27,1).Value = "season_spring"
$ws.Cells.Item(27,2).Value = "Spring"
$ws.Cells.Item(28,1).Value = "season_summer"
$ws.Cells.Item(28,2).Value = "Summer"
$ws.Cells.Item(29,1).Value = "season_autumn"
$ws.Cells.Item(29,2).Value = "Autumn"
$ws.Cells.Item(30,1).Value = "atmosphere_altitude"
$ws.Cells.Item(30,2).Value = "Altitude"
$ws.Cells.Item(31,1).Value = "atmosphere_humidity"
$ws.Cells.Item(31,2).Value = "Humidity"
$ws.Cells.Item(32,1).Value = "atmosphere_temperature"
$ws.Cells.Item(32,2).Value = "Temperature"
$ws.Cells.Item(33,1).Value = "atmosphere_windStrength"
$ws.Cells.Item(33,2).Value = "Wind Strength"
$ws.Cells.Item(34,1).Value = "climate_title"
$ws.Cells.Item(34,2).Value = "Climate"
$ws.Cells.Item(35,1).Value = "climate_temperate"
$ws.Cells.Item(35,2).Value = "Temperate"
$ws.Cells.Item(36,1).Value = "climate_tropical"
$ws.Cells.Item(36,2).Value = "Tropical"
$ws.Cells.Item(37,1).Value = "climate_oceanic"
$ws.Cells.Item(37,2).Value = "Oceanic"
$ws.Cells.Item(38,1).Value = "climate_desert"
$ws.Cells.Item(38,2).Value = "Desert"
$ws.Cells.Item(39,1).Value = "climate_tundra"
$ws.Cells.Item(39,2).Value = "Tundra"
$ws.Cells.Item(40,1).Value = "climate_highland"
$ws.Cells.Item(40,2).Value = "Highland"
$ws.Cells.Item(41,1).Value = "climate_mediterranean"
$ws.Cells.Item(41,2).Value = "Mediterranean"
$ws.Cells.Item(42,1).Value = "region_title"
$ws.Cells.Item(42,2).Value = "Region"
$ws.Cells.Item(43,1).Value = "region_NA"
$ws.Cells.Item(43,2).Value = "North American Great Plains"
$ws.Cells.Item(44,1).Value = "region_PH"
$ws.Cells.Item(44,2).Value = "Northern Luzon, Philippines"
$ws.Cells.Item(45,1).Value = "region_GB"
$ws.Cells.Item(45,2).Value = "British Isles"
$ws.Cells.Item(46,1).Value = "region_MG"
$ws.Cells.Item(46,2).Value = "Madagascar"
$ws.Cells.Item(47,1).Value = "region_EG"
$ws.Cells.Item(47,2).Value = "Egypt"
$ws.Cells.Item(48,1).Value = "region_GL"
$ws.Cells.Item(48,2).Value = "Greenland"
$ws.Cells.Item(49,1).Value = "region_BR"
$ws.Cells.Item(49,2).Value = "Brazil"
$ws.Cells.Item(50,1).Value = "region_CL"
$ws.Cells.Item(50,2).Value = "Andes Mountains, Chile"
$ws.Cells.Item(51,1).Value = "region_IT"
$ws.Cells.Item(51,2).Value = "Italy"
$ws.Cells.Item(52,1).Value = "region_AU"
$ws.Cells.Item(52,2).Value = "Australian Outback"
$ws.Cells.Item(53,1).Value = "weatherForecast"
$ws.Cells.Item(53,2).Value = "Weather Forecast"
$ws.Cells.Item(54,1).Value = "day_today"
$ws.Cells.Item(54,2).Value = "Today"
$ws.Cells.Item(55,1).Value = "day_monday"
$ws.Cells.Item(55,2).Value = "Monday"
$ws.Cells.Item(56,1).Value = "day_tuesday"
$ws.Cells.Item(56,2).Value = "Tuesday"
$ws.Cells.Item(57,1).Value = "day_wednesday"
$ws.Cells.Item(57,2).Value = "Wednesday"
$ws.Cells.Item(58,1).Value = "day_thursday"
$ws.Cells.Item(58,2).Value = "Thursday"
$ws.Cells.Item(59,1).Value = "day_friday"
$ws.Cells.Item(59,2).Value = "Friday"
$ws.Cells.Item(60,1).Value = "day_saturday"
$ws.Cells.Item(60,2).Value = "Saturday"
$ws.Cells.Item(61,1).Value = "day_sunday"
$ws.Cells.Item(61,2).Value = "Sunday"
$ws.Cells.Item(62,1).Value = "weather_sunny"
$ws.Cells.Item(62,2).Value = "Sunny"
$ws.Cells.Item(63,1).Value = "weather_partly_sunny"
$ws.Cells.Item(63,2).Value = "Partly Sunny"
$ws.Cells.Item(64,1).Value = "weather_mostly_cloudy"
$ws.Cells.Item(64,2).Value = "Mostly Cloudy"
$ws.Cells.Item(65,1).Value = "weather_cloudy"
$ws.Cells.Item(65,2).Value = "Cloudy"
$ws.Cells.Item(66,1).Value = "weather_light_rain"
$ws.Cells.Item(66,2).Value = "Light Rain"
$ws.Cells.Item(67,1).Value = "weather_rain"
$ws.Cells.Item(67,2).Value = "Rain"
$ws.Cells.Item(68,1).Value = "weather_heavy_rain"
$ws.Cells.Item(68,2).Value = "Heavy Rain"
$ws.Cells.Item(69,1).Value = "weather_typhoon"
$ws.Cells.Item(69,2).Value = "Typhoon"
$ws.Cells.Item(70,1).Value = "weather_haze"
$ws.Cells.Item(70,2).Value = "Haze"
$ws.Cells.Item(71,1).Value = "weather_dustStorm"
$ws.Cells.Item(71,2).Value = "Dust Storm"
$ws.Cells.Item(72,1).Value = "weather_light_snow"
$ws.Cells.Item(72,2).Value = "Light Snow"
$ws.Cells.Item(72,2).WrapText = $true
$ws.Cells.Item(73,1).Value = "weather_snow"
$ws.Cells.Item(73,2).Value = "Snow"
$ws.Cells.Item(73,2).WrapText = $true

# New rows 74-79: intro dialogue keys (no special style on column B)
$ws.Cells.Item(74,1).Value = "intro_0_0"
$ws.Cells.Item(74,2).Value = "Unidentified ships approaching Earth!"
$ws.Cells.Item(75,1).Value = "intro_1_0"
$ws.Cells.Item(75,2).Value = "They appear to be frog-like. Let me put on my frog suit to communicate with these peculiar creatures!"
$ws.Cells.Item(76,1).Value = "intro_2_0"
$ws.Cells.Item(76,2).Value = "That’s better! Let’s see if they are of any threat…"
$ws.Cells.Item(77,1).Value = "intro_3_0"
$ws.Cells.Item(77,2).Value = "It looks like they have been exiled from their planet, and are looking for a new home."
$ws.Cells.Item(78,1).Value = "intro_3_1"
$ws.Cells.Item(78,2).Value = "Well, we can’t just let them hang about in space. Besides, how often are we visited by sentient aliens from outer space?"
$ws.Cells.Item(79,1).Value = "intro_4_0"
$ws.Cells.Item(79,2).Value = "They are expressing their gratitude, and are ready to cooperate in exchange for their refuge."

# New rows 80-81: intro dialogue keys with vertical-center alignment on column B
$ws.Cells.Item(80,1).Value = "intro_4_1"
$ws.Cells.Item(80,2).Value = "Well, why not? Let’s give these hapless frogs some proper homes to settle in. There’s still plenty of room here on Earth."
$ws.Cells.Item(80,2).VerticalAlignment = -4108
$ws.Cells.Item(81,1).Value = "intro_4_2"
$ws.Cells.Item(81,2).Value = "Now commencing operation: Project Bloom – A quest for home!"
$ws.Cells.Item(81,2).VerticalAlignment = -4108

# Update selection/active cell to match final view state
$ws.Range("B80").Select()

